$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04763786555579896
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 10.53319022774986

$ws.Range("B3").Value = 3.230985683306322
$ws.Range("C3").Value = 10.29869402782916
$ws.Range("D3").Value = 0.1575252929769615
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 22.34743749006142

$ws.Range("B4").Value = 3.230985683306322
$ws.Range("C4").Value = 1.667794583268128
$ws.Range("D4").Value = 0.1575252929769615
$ws.Range("E4").Value = 0.496779210170732
$ws.Range("G4").Value = 5.553084769722144
